$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new row of hospitalization data for 4 April 2020 (NY State update)
# Clone formatting from the prior row (20) so the new row matches the
# existing date/number styling exactly, then fill in the new values.
$ws.Range("A20:E20").Copy()
$ws.Range("A21:E21").PasteSpecial(-4122)

$ws.Range("A21").Value = 43925
$ws.Range("B21").Value = 574
$ws.Range("C21").Value = 250
$ws.Range("D21").Value = 1709
$ws.Range("E21").Value = 316

# Update selection to mirror the authored workbook state
$ws.Range("D22").Select()
